# Unit 5 Assignment Brief: swap the logical "name" Word shows for the two
# header/footer logo pictures (BTec_Logo-Orange <-> PearsonLogo) across
# every section's headers and footers.
#
# image1.png -> image2.png   (Pearson Edexcel logo, in the footers)
# image2.jpg -> image1.jpg   (BTec logo, in the headers)

$d = $word.ActiveDocument

function Rename-LogoInlineShapes($range) {
    if ($null -eq $range) { return }
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image2.png"
        } elseif ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image1.jpg"
        }
    }
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $section = $d.Sections.Item($si)

    for ($hi = 1; $hi -le $section.Headers.Count; $hi++) {
        $header = $section.Headers.Item($hi)
        if ($header.Exists) {
            Rename-LogoInlineShapes $header.Range
        }
    }

    for ($fi = 1; $fi -le $section.Footers.Count; $fi++) {
        $footer = $section.Footers.Item($fi)
        if ($footer.Exists) {
            Rename-LogoInlineShapes $footer.Range
        }
    }
}
